# Regenerate save_data column G ("K", strikeouts) using the new source
# stat (K) instead of the old "Strike#" derived value. Only column G
# values change; row 12 keeps its original value (0 -> 0, unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    13 = 1
    14 = 2
    15 = 2
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
